$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 386, shifting the existing rows 386:411 down to 387:412
$ws.Rows.Item(386).Insert()

# Populate the newly inserted row 386 with the new record
$ws.Range("A386").Value = 9
$ws.Range("B386").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C386").Value = "Metropolitana"
$ws.Range("D386").Value = 45021
$ws.Range("D386").NumberFormat = $ws.Range("D387").NumberFormat
$ws.Range("E386").Value = 13
$ws.Range("F386").Value = 100112043
$ws.Range("G386").Value = "Pepino ensalada"
$ws.Range("H386").Value = "Sin especificar"
$ws.Range("I386").Value = "Primera"
$ws.Range("J386").Value = 65
$ws.Range("K386").Value = 8000
$ws.Range("L386").Value = 9000
$ws.Range("M386").Value = 8538
$ws.Range("N386").Value = "$/caja 60 unidades"
$ws.Range("O386").Value = "Región de Arica y Parinacota"
$ws.Range("P386").Value = 142
$ws.Range("Q386").Value = 60
$ws.Range("R386").Value = "Hortaliza"
